$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 4, pushing the existing rows 4-5 down to 6-7.
$ws.Range("A4:A5").EntireRow.Insert()

# New row 4 data
$ws.Range("A4").Value = "ruMjnUSk"
$ws.Range("B4").Value = "24/03/2025"
$ws.Range("C4").Value = "16:00"
$ws.Range("D4").Value = "COLOMBIA - PRIMERA B"
$ws.Range("E4").Value = "Bogota"
$ws.Range("F4").Value = "Patriotas"
$ws.Range("G4").Value = 4.5
$ws.Range("H4").Value = 3.35
$ws.Range("I4").Value = 1.75
$ws.Range("J4").Value = 4.8
$ws.Range("K4").Value = 2.05
$ws.Range("L4").Value = 2.37
$ws.Range("O4").Value = 1.32
$ws.Range("P4").Value = 2.85
$ws.Range("S4").Value = 1.93
$ws.Range("T4").Value = 1.7
$ws.Range("W4").Value = 3.15
$ws.Range("X4").Value = 1.26
$ws.Range("Y4").Value = 1.42
$ws.Range("Z4").Value = 2.45
$ws.Range("AA4").Value = 1.83
$ws.Range("AB4").Value = 1.78
$ws.Range("AC4").Value = 11.75
$ws.Range("AD4").Value = 25
$ws.Range("AE4").Value = 14.5
$ws.Range("AF4").Value = 80
$ws.Range("AG4").Value = 45
$ws.Range("AH4").Value = 50
$ws.Range("AI4").Value = 8.75
$ws.Range("AJ4").Value = 6.6
$ws.Range("AK4").Value = 16
$ws.Range("AL4").Value = 80
$ws.Range("AM4").Value = 6.3
$ws.Range("AN4").Value = 7.8
$ws.Range("AO4").Value = 8.25
$ws.Range("AP4").Value = 14
$ws.Range("AQ4").Value = 14.5
$ws.Range("AR4").Value = 29
$ws.Range("AS4").Value = 700

# New row 5 data
$ws.Range("A5").Value = "WvRAeY53"
$ws.Range("B5").Value = "24/03/2025"
$ws.Range("C5").Value = "22:20"
$ws.Range("D5").Value = "COLOMBIA - PRIMERA B"
$ws.Range("E5").Value = "Leones"
$ws.Range("F5").Value = "Real Cundinamarca"
$ws.Range("G5").Value = 2.3
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 2.87
$ws.Range("J5").Value = 2.9
$ws.Range("K5").Value = 2.05
$ws.Range("L5").Value = 3.45
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 3
$ws.Range("S5").Value = 1.85
$ws.Range("T5").Value = 1.75
$ws.Range("W5").Value = 2.95
$ws.Range("X5").Value = 1.3
$ws.Range("Y5").Value = 1.4
$ws.Range("Z5").Value = 2.52
$ws.Range("AA5").Value = 1.65
$ws.Range("AB5").Value = 1.98
$ws.Range("AC5").Value = 8.25
$ws.Range("AD5").Value = 11.75
$ws.Range("AE5").Value = 9
$ws.Range("AF5").Value = 24
$ws.Range("AG5").Value = 18.5
$ws.Range("AH5").Value = 27
$ws.Range("AI5").Value = 9.75
$ws.Range("AJ5").Value = 6.3
$ws.Range("AK5").Value = 13.5
$ws.Range("AL5").Value = 60
$ws.Range("AM5").Value = 9
$ws.Range("AN5").Value = 15
$ws.Range("AO5").Value = 10.5
$ws.Range("AP5").Value = 35
$ws.Range("AQ5").Value = 25
$ws.Range("AR5").Value = 32
$ws.Range("AS5").Value = 450
